$wb = $excel.ActiveWorkbook

# Remember which sheet/cell was originally active so the selection can be
# restored at the end (adding a sheet shouldn't change the user's view).
$originalActiveSheet = $wb.Worksheets.Item(1)

# Add the new worksheet after the last existing sheet so it becomes the
# 4th tab, named "ODI Batting Extra".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Match the page margins / outline defaults used by the other sheets in
# this workbook.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Force the whole used range to text format so the string-looking values
# ("4187", "0", "1.60%", ...) are stored as text instead of being
# auto-coerced into numbers/percentages.
$newSheet.Range("A1:F3").NumberFormat = "@"

# Header row - styled like the header rows on the other sheets
# (bold font, thin box border, centered horizontal / top vertical alignment)
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Row 2
$newSheet.Cells.Item(2, 1).Value = "4187"
$newSheet.Cells.Item(2, 3).Value = "0"
$newSheet.Cells.Item(2, 4).Value = "0"
$newSheet.Cells.Item(2, 5).Value = "1.60%"
$newSheet.Cells.Item(2, 6).Value = "NO"

# BATTING_POSITION is a genuine number, not text
$newSheet.Range("B2").NumberFormat = "General"
$newSheet.Cells.Item(2, 2).Value = 10

# Row 3 - only MATCH_CODE and MAN_OF_MATCH are populated
$newSheet.Cells.Item(3, 1).Value = "4188"
$newSheet.Cells.Item(3, 6).Value = "NO"

# Restore the original selection/active sheet
$originalActiveSheet.Activate()
$originalActiveSheet.Range("A1").Select()
